$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 16366
$ws.Cells.Item(18, 9).Value = 549.5
$ws.Cells.Item(18, 11).Value = 549.5
$ws.Cells.Item(18, 13).Value = -265.5
$ws.Cells.Item(19, 8).Value = 3699.2
$ws.Cells.Item(19, 9).Value = 3498
$ws.Cells.Item(19, 10).Value = 4001
$ws.Cells.Item(19, 11).Value = 3498
$ws.Cells.Item(19, 12).Value = 4001
$ws.Cells.Item(19, 13).Value = -3323
$ws.Cells.Item(19, 14).Value = -4351
$ws.Cells.Item(33, 8).Value = 691466.1
$ws.Cells.Item(33, 9).Value = 1150650.9
$ws.Cells.Item(33, 10).Value = 2688.9
$ws.Cells.Item(33, 11).Value = 1150650.9
$ws.Cells.Item(33, 12).Value = 2688.9
$ws.Cells.Item(33, 13).Value = -1150421.9
$ws.Cells.Item(33, 14).Value = -3146.9
$ws.Cells.Item(40, 8).Value = 3112.4167
$ws.Cells.Item(40, 9).Value = 1287.5
$ws.Cells.Item(40, 11).Value = 1287.5
$ws.Cells.Item(40, 13).Value = -1112.5
$ws.Cells.Item(43, 8).Value = 5832.8335
$ws.Cells.Item(43, 9).Value = 3500
$ws.Cells.Item(43, 11).Value = 3500
$ws.Cells.Item(43, 13).Value = -3431
$ws.Cells.Item(76, 8).Value = 9999.666999999999
$ws.Cells.Item(76, 9).Value = 10000
$ws.Cells.Item(76, 10).Value = 9999
$ws.Cells.Item(76, 11).Value = 10000
$ws.Cells.Item(76, 12).Value = 9999
$ws.Cells.Item(76, 13).Value = -9685
$ws.Cells.Item(76, 14).Value = -10629
$ws.Cells.Item(79, 8).Value = 9999.666999999999
$ws.Cells.Item(79, 9).Value = 10000
$ws.Cells.Item(79, 10).Value = 9999
$ws.Cells.Item(79, 11).Value = 10000
$ws.Cells.Item(79, 12).Value = 9999
$ws.Cells.Item(79, 13).Value = -8908
$ws.Cells.Item(79, 14).Value = -12183
$ws.Cells.Item(94, 8).Value = 970.6667
$ws.Cells.Item(94, 9).Value = 971.5
$ws.Cells.Item(94, 11).Value = 971.5
$ws.Cells.Item(94, 13).Value = -520.5
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:L134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:L139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7686147
$ws.Cells.Item(32, 9).Value = 1075284.1
$ws.Cells.Item(32, 11).Value = 1075284.1
$ws.Cells.Item(32, 13).Value = -1074997.1
$ws.Cells.Item(45, 8).Value = 15575.174
$ws.Cells.Item(45, 9).Value = 13779.177
$ws.Cells.Item(45, 10).Value = 20663.834
$ws.Cells.Item(45, 11).Value = 13779.177
$ws.Cells.Item(45, 12).Value = 20663.834
$ws.Cells.Item(45, 13).Value = -13402.177
$ws.Cells.Item(45, 14).Value = -21417.834
$ws.Cells.Item(61, 8).Value = 2612.0833
$ws.Cells.Item(61, 9).Value = 2513.5557
$ws.Cells.Item(61, 11).Value = 2513.5557
$ws.Cells.Item(61, 13).Value = -2301.5557
$ws.Cells.Item(74, 8).Value = 2101
$ws.Cells.Item(74, 9).Value = 1755.5
$ws.Cells.Item(74, 11).Value = 1755.5
$ws.Cells.Item(74, 13).Value = -881.5
$ws.Cells.Item(77, 8).Value = 2101
$ws.Cells.Item(77, 9).Value = 1755.5
$ws.Cells.Item(77, 11).Value = 8777.5
$ws.Cells.Item(77, 13).Value = -4409.5
$ws.Cells.Item(110, 8).Value = 1497.5769
$ws.Cells.Item(110, 9).Value = 1061.1052
$ws.Cells.Item(110, 11).Value = 1061.1052
$ws.Cells.Item(110, 13).Value = 983.8948
$ws.Cells.Item(122, 8).Value = 10823
$ws.Cells.Item(122, 9).Value = 15219.177
$ws.Cells.Item(122, 10).Value = 3349.5
$ws.Cells.Item(122, 11).Value = 45657.531
$ws.Cells.Item(122, 12).Value = 10048.5
$ws.Cells.Item(122, 13).Value = -43207.531
$ws.Cells.Item(122, 14).Value = -14948.5
$ws.Cells.Item(132, 8).Value = 2446.0159
$ws.Cells.Item(132, 9).Value = 2159.9443
$ws.Cells.Item(132, 11).Value = 6479.8329
$ws.Cells.Item(132, 13).Value = -3949.8329
$ws.Cells.Item(136, 8).Value = 2612.0833
$ws.Cells.Item(136, 9).Value = 2513.5557
$ws.Cells.Item(136, 11).Value = 7540.6671
$ws.Cells.Item(136, 13).Value = -4990.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4194.8237
$ws.Cells.Item(20, 9).Value = 2700.3076
$ws.Cells.Item(20, 10).Value = 9052
$ws.Cells.Item(20, 11).Value = 2700.3076
$ws.Cells.Item(20, 12).Value = 9052
$ws.Cells.Item(20, 13).Value = -2453.3076
$ws.Cells.Item(20, 14).Value = -9546
$ws.Cells.Item(107, 8).Value = 5750
$ws.Cells.Item(107, 9).Value = 5666.6665
$ws.Cells.Item(107, 11).Value = 5666.6665
$ws.Cells.Item(107, 13).Value = -3746.6665
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(120, 8).Value = 80000
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 80000
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 80000
$ws.Cells.Item(120, 14).Value = -89676
$ws.Cells.Item(122, 8).Value = 60000
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 60000
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 60000
$ws.Cells.Item(122, 14).Value = -69800
$ws.Cells.Item(123, 8).Value = 70000
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 70000
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 70000
$ws.Cells.Item(123, 14).Value = -79800
$ws.Cells.Item(124, 8).Value = 120000
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 120000
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 120000
$ws.Cells.Item(124, 14).Value = -129820
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(126, 8).Value = 70774.5
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 70774.5
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 70774.5
$ws.Cells.Item(126, 14).Value = -80654.5
$ws.Cells.Item(127, 8).Value = 60780
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 60780
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 60780
$ws.Cells.Item(127, 14).Value = -70700
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 40000
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 40000
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 40000
$ws.Cells.Item(129, 14).Value = -50000
$ws.Cells.Item(130, 8).Value = 100780
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 100780
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 100780
$ws.Cells.Item(130, 14).Value = -110820
$ws.Cells.Item(131, 8).Value = 40000
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 40000
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 40000
$ws.Cells.Item(131, 14).Value = -50080
$ws.Cells.Item(132, 8).Value = 176543
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 176543
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 176543
$ws.Cells.Item(132, 14).Value = -186663
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(134, 8).Value = 4532.116
$ws.Cells.Item(134, 9).Value = 4259.1353
$ws.Cells.Item(134, 10).Value = 6215.5
$ws.Cells.Item(134, 11).Value = 12777.4059
$ws.Cells.Item(134, 12).Value = 18646.5
$ws.Cells.Item(134, 13).Value = -10242.4059
$ws.Cells.Item(134, 14).Value = -23716.5
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 5000
$ws.Cells.Item(2, 10).Value = 5000
$ws.Cells.Item(2, 12).Value = 5000
$ws.Cells.Item(2, 14).Value = -5226
$ws.Cells.Item(51, 8).Value = 12142.857
$ws.Cells.Item(51, 10).Value = 15000
$ws.Cells.Item(51, 12).Value = 15000
$ws.Cells.Item(51, 14).Value = -16472
$ws.Cells.Item(58, 8).Value = 1364
$ws.Cells.Item(58, 9).Value = 1330.4
$ws.Cells.Item(58, 11).Value = 1330.4
$ws.Cells.Item(58, 13).Value = -1127.4
$ws.Cells.Item(61, 8).Value = 12142.857
$ws.Cells.Item(61, 10).Value = 15000
$ws.Cells.Item(61, 12).Value = 15000
$ws.Cells.Item(61, 14).Value = -15696
$ws.Cells.Item(68, 8).Value = 25000
$ws.Cells.Item(68, 10).Value = 25000
$ws.Cells.Item(68, 12).Value = 25000
$ws.Cells.Item(68, 14).Value = -26498
$ws.Cells.Item(70, 8).Value = 23090
$ws.Cells.Item(70, 10).Value = 23090
$ws.Cells.Item(70, 12).Value = 23090
$ws.Cells.Item(70, 14).Value = -23720
$ws.Cells.Item(71, 8).Value = 25000
$ws.Cells.Item(71, 10).Value = 25000
$ws.Cells.Item(71, 12).Value = 75000
$ws.Cells.Item(71, 14).Value = -82488
$ws.Cells.Item(73, 8).Value = 23090
$ws.Cells.Item(73, 10).Value = 23090
$ws.Cells.Item(73, 12).Value = 23090
$ws.Cells.Item(73, 14).Value = -25274
$ws.Cells.Item(86, 8).Value = 41671304
$ws.Cells.Item(86, 9).Value = 90913020
$ws.Cells.Item(86, 11).Value = 90913020
$ws.Cells.Item(86, 13).Value = -90911897
$ws.Cells.Item(89, 8).Value = 41671304
$ws.Cells.Item(89, 9).Value = 90913020
$ws.Cells.Item(89, 11).Value = 454565100
$ws.Cells.Item(89, 13).Value = -454559484
$ws.Cells.Item(107, 8).Value = 2112.5789
$ws.Cells.Item(107, 9).Value = 2209
$ws.Cells.Item(107, 10).Value = 1980
$ws.Cells.Item(107, 11).Value = 2209
$ws.Cells.Item(107, 12).Value = 1980
$ws.Cells.Item(107, 13).Value = -289
$ws.Cells.Item(107, 14).Value = -5820
$ws.Cells.Item(132, 8).Value = 2351.2856
$ws.Cells.Item(132, 9).Value = 2258.739
$ws.Cells.Item(132, 10).Value = 2777
$ws.Cells.Item(132, 11).Value = 6776.217000000001
$ws.Cells.Item(132, 12).Value = 8331
$ws.Cells.Item(132, 13).Value = -4246.217000000001
$ws.Cells.Item(132, 14).Value = -13391
$ws.Cells.Item(134, 8).Value = 2663.1
$ws.Cells.Item(134, 9).Value = 2080.2354
$ws.Cells.Item(134, 10).Value = 5966
$ws.Cells.Item(134, 11).Value = 6240.706200000001
$ws.Cells.Item(134, 12).Value = 17898
$ws.Cells.Item(134, 13).Value = -3705.706200000001
$ws.Cells.Item(134, 14).Value = -22968
$ws.Cells.Item(136, 8).Value = 1364
$ws.Cells.Item(136, 9).Value = 1330.4
$ws.Cells.Item(136, 11).Value = 3991.2
$ws.Cells.Item(136, 13).Value = -1441.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 14273937
$ws.Cells.Item(4, 9).Value = 3342585.8
$ws.Cells.Item(4, 10).Value = 44881720
$ws.Cells.Item(4, 11).Value = 10027757.4
$ws.Cells.Item(4, 12).Value = 134645160
$ws.Cells.Item(4, 13).Value = -10027645.4
$ws.Cells.Item(4, 14).Value = -134645384
$ws.Cells.Item(80, 8).Value = 17642856
$ws.Cells.Item(80, 9).Value = 5000
$ws.Cells.Item(80, 10).Value = 23522140
$ws.Cells.Item(80, 11).Value = 15000
$ws.Cells.Item(80, 12).Value = 70566420
$ws.Cells.Item(80, 13).Value = -14064
$ws.Cells.Item(80, 14).Value = -70568292
$ws.Cells.Item(83, 8).Value = 17642856
$ws.Cells.Item(83, 9).Value = 5000
$ws.Cells.Item(83, 10).Value = 23522140
$ws.Cells.Item(83, 11).Value = 45000
$ws.Cells.Item(83, 12).Value = 211699260
$ws.Cells.Item(83, 13).Value = -40320
$ws.Cells.Item(83, 14).Value = -211708620
$ws.Cells.Item(93, 8).Value = 129909.875
$ws.Cells.Item(93, 10).Value = 5611.4287
$ws.Cells.Item(93, 12).Value = 16834.2861
$ws.Cells.Item(93, 14).Value = -20578.2861
$ws.Cells.Item(107, 8).Value = 1034.1818
$ws.Cells.Item(107, 10).Value = 1152.0555
$ws.Cells.Item(107, 12).Value = 3456.1665
$ws.Cells.Item(107, 14).Value = -7296.166499999999
$ws.Cells.Item(132, 8).Value = 1816.909
$ws.Cells.Item(132, 9).Value = 1713
$ws.Cells.Item(132, 11).Value = 15417
$ws.Cells.Item(132, 13).Value = -12887

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 15240.857
$ws.Cells.Item(46, 9).Value = 1339.2
$ws.Cells.Item(46, 11).Value = 1339.2
$ws.Cells.Item(46, 13).Value = -1183.2
$ws.Cells.Item(70, 8).Value = 6667.6562
$ws.Cells.Item(70, 10).Value = 8416.333000000001
$ws.Cells.Item(70, 12).Value = 8416.333000000001
$ws.Cells.Item(70, 14).Value = -8956.333000000001
$ws.Cells.Item(73, 8).Value = 6667.6562
$ws.Cells.Item(73, 10).Value = 8416.333000000001
$ws.Cells.Item(73, 12).Value = 8416.333000000001
$ws.Cells.Item(73, 14).Value = -10288.333
$ws.Cells.Item(113, 8).Value = 13728.8
$ws.Cells.Item(113, 9).Value = 19548.834
$ws.Cells.Item(113, 10).Value = 4998.75
$ws.Cells.Item(113, 11).Value = 19548.834
$ws.Cells.Item(113, 12).Value = 4998.75
$ws.Cells.Item(113, 13).Value = -17378.834
$ws.Cells.Item(113, 14).Value = -9338.75
$ws.Cells.Item(122, 8).Value = 3009
$ws.Cells.Item(122, 9).Value = 3009
$ws.Cells.Item(122, 11).Value = 9027
$ws.Cells.Item(122, 13).Value = -6577
$ws.Cells.Item(132, 8).Value = 3499.7334
$ws.Cells.Item(132, 9).Value = 3519.7778
$ws.Cells.Item(132, 10).Value = 3419.5557
$ws.Cells.Item(132, 11).Value = 10559.3334
$ws.Cells.Item(132, 12).Value = 10258.6671
$ws.Cells.Item(132, 13).Value = -8029.3334
$ws.Cells.Item(132, 14).Value = -15318.6671
$ws.Cells.Item(135, 8).Value = 74999.82000000001
$ws.Cells.Item(135, 10).Value = 74999.82000000001
$ws.Cells.Item(135, 12).Value = 74999.82000000001
$ws.Cells.Item(135, 14).Value = -85139.82000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1295
$ws.Cells.Item(93, 10).Value = 2628.3333
$ws.Cells.Item(93, 12).Value = 2628.3333
$ws.Cells.Item(93, 14).Value = -5124.3333
$ws.Cells.Item(122, 8).Value = 14402.3
$ws.Cells.Item(122, 9).Value = 16159.375
$ws.Cells.Item(122, 10).Value = 7374
$ws.Cells.Item(122, 11).Value = 48478.125
$ws.Cells.Item(122, 12).Value = 22122
$ws.Cells.Item(122, 13).Value = -46028.125
$ws.Cells.Item(122, 14).Value = -27022
$ws.Cells.Item(132, 8).Value = 6017.4863
$ws.Cells.Item(132, 9).Value = 4197.72
$ws.Cells.Item(132, 10).Value = 9808.666999999999
$ws.Cells.Item(132, 11).Value = 12593.16
$ws.Cells.Item(132, 12).Value = 29426.001
$ws.Cells.Item(132, 13).Value = -10063.16
$ws.Cells.Item(132, 14).Value = -34486.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1334000
$ws.Cells.Item(81, 9).Value = 1400800
$ws.Cells.Item(81, 11).Value = 2801600
$ws.Cells.Item(81, 13).Value = -2800539
$ws.Cells.Item(84, 8).Value = 1334000
$ws.Cells.Item(84, 9).Value = 1400800
$ws.Cells.Item(84, 11).Value = 14008000
$ws.Cells.Item(84, 13).Value = -14002696
$ws.Cells.Item(100, 8).Value = 3297.7693
$ws.Cells.Item(100, 9).Value = 5003
$ws.Cells.Item(100, 11).Value = 10006
$ws.Cells.Item(100, 13).Value = -9465
$ws.Cells.Item(122, 8).Value = 1851
$ws.Cells.Item(122, 9).Value = 1200
$ws.Cells.Item(122, 10).Value = 2013.75
$ws.Cells.Item(122, 11).Value = 3600
$ws.Cells.Item(122, 12).Value = 6041.25
$ws.Cells.Item(122, 13).Value = -1150
$ws.Cells.Item(122, 14).Value = -10941.25
$ws.Cells.Item(132, 8).Value = 11008.6875
$ws.Cells.Item(132, 9).Value = 12366.786
$ws.Cells.Item(132, 10).Value = 1502
$ws.Cells.Item(132, 11).Value = 37100.358
$ws.Cells.Item(132, 12).Value = 4506
$ws.Cells.Item(132, 13).Value = -34570.358
$ws.Cells.Item(132, 14).Value = -9566
